$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append new workout-log rows (31-43) to the "Activity" sheet, continuing the
# existing table (dates 31 Oct 2023 - 9 Nov 2023), and bring in the handful
# of new "Details" labels that come with them.
# ---------------------------------------------------------------------------

# Reuse number formats already present on the sheet so no new cell styles
# get created in styles.xml:
#   Row 2 -> A:G "activity with distance" pattern (s2,s1,s1,s1,s1,s3,s1)
#   Row 4 -> A:D,F:G "activity without distance" pattern (s2,s1,s1,s1,s3,s1)
#   Row 3 -> A:B "Rest" pattern (s2,s1)

function New-FullRow($row, $date, $activity, $details, $duration, $distance, $time, $intensity) {
    $ws.Range("A2:G2").Copy($ws.Range("A${row}:G${row}"))
    $ws.Range("A$row").Value = $date
    $ws.Range("B$row").Value = $activity
    $ws.Range("C$row").Value = $details
    $ws.Range("D$row").Value = $duration
    $ws.Range("E$row").Value = $distance
    $ws.Range("F$row").Value = $time
    $ws.Range("G$row").Value = $intensity
}

function New-NoDistanceRow($row, $date, $activity, $details, $duration, $time, $intensity) {
    $ws.Range("A4:D4").Copy($ws.Range("A${row}:D${row}"))
    $ws.Range("F4:G4").Copy($ws.Range("F${row}:G${row}"))
    $ws.Range("A$row").Value = $date
    $ws.Range("B$row").Value = $activity
    $ws.Range("C$row").Value = $details
    $ws.Range("D$row").Value = $duration
    $ws.Range("F$row").Value = $time
    $ws.Range("G$row").Value = $intensity
}

function New-RestRow($row, $date, $activity) {
    $ws.Range("A3:B3").Copy($ws.Range("A${row}:B${row}"))
    $ws.Range("A$row").Value = $date
    $ws.Range("B$row").Value = $activity
}

# Row 31 - Swimming
New-FullRow 31 45230 "Swimming" "Breatstroke / Frontcrawl" 30 525 0.34027777777777773 5

# Row 32 - Gym (Back)
New-NoDistanceRow 32 45230 "Gym" "Back" 90 0.73958333333333337 6

# Row 33 - Rest
New-RestRow 33 45231 "Rest"

# Row 34 - Swimming
New-FullRow 34 45232 "Swimming" "Breaststroke / Frontcrawl" 35 800 0.34722222222222227 6

# Row 35 - Rest
New-RestRow 35 45233 "Rest"

# Row 36 - Running ( Short Run) -- first brand-new shared string
New-FullRow 36 45234 "Running" " Short Run" 45 7150 0.4069444444444445 6

# Row 37 - Gym (Legs)
New-NoDistanceRow 37 45234 "Gym" "Legs" 120 0.71875 8

# Row 38 - Gym (Shoulders / Triceps) -- second brand-new shared string
New-NoDistanceRow 38 45235 "Gym" "Shoulders / Triceps" 90 0.65555555555555556 7

# Row 39 - Swimming
New-FullRow 39 45236 "Swimming" "Breastroke / Frontcrawl" 35 750 0.33333333333333331 8

# Row 41 is filled in before row 40 is finished off below, so the new
# shared-string table ends up in the same order the workbook author
# produced it in (Core / Chest, then Core, then Abs).
# Row 41 - Gym (Core / Chest)
New-NoDistanceRow 41 45237 "Gym" "Core / Chest" 75 0.71875 7

# Row 40 - Core (Abs)
New-NoDistanceRow 40 45236 "Core" "Abs" 10 0.79166666666666663 6

# Row 42 - Swimming
New-FullRow 42 45238 "Swimming" "Breastroke / Frontcrawl" 30 600 0.33333333333333331 5

# Row 43 - Rest
New-RestRow 43 45239 "Rest"

# Scroll / selection state, mirroring where the author ended up after typing.
$ws.Range("D44").Select()
